$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('T2').Value = 'maa://22742 (91.5), *maa://20791 (62.86)'
$ws.Range('AF2').Value = 'maa://25251 (91.49), ***maa://21730 (19.4), ***maa://39501 (19.05), *maa://36675 (60.0)'
$ws.Range('H3').Value = 'maa://21247 (98.4), *maa://22748 (60.0)'
$ws.Range('L3').Value = '*maa://22880 (65.93), maa://20276 (84.71), *maa://22749 (72.73)'
$ws.Range('AB3').Value = 'maa://24390 (96.55)'
$ws.Range('T4').Value = 'maa://32509 (97.94), maa://27295 (83.61), maa://22754 (91.67), *maa://21746 (55.81), *maa://31008 (78.05)'
$ws.Range('X4').Value = '**maa://32495 (47.69), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (89.66)'
$ws.Range('D5').Value = 'maa://21245 (83.1), maa://22744 (84.0)'
$ws.Range('AB5').Value = '*maa://29863 (71.88), ***maa://22752 (13.33), **maa://26013 (42.86)'
$ws.Range('D6').Value = 'maa://42407 (94.44)'
$ws.Range('AF7').Value = '*maa://26191 (67.95), *maa://36671 (70.83), *maa://42530 (69.23)'
$ws.Range('A8').Value = '更新日期：2024.12.18 13:18:43'
$ws.Range('X8').Value = 'maa://21411 (95.89)'
$ws.Range('AF8').Value = '*maa://24479 (77.78), *maa://21990 (53.85)'
$ws.Range('AF9').Value = 'maa://26206 (89.32), **maa://22865 (50.0)'
$ws.Range('D10').Value = '***maa://25695 (19.44), **maa://32237 (40.48), ***maa://34206 (18.18), ***maa://39951 (17.07), ***maa://39243 (28.57)'
$ws.Range('T11').Value = 'maa://22747 (93.42), maa://22501 (98.55)'
$ws.Range('X11').Value = 'maa://36713 (98.15)'
$ws.Range('AB12').Value = 'maa://23669 (95.36), maa://36677 (93.88), maa://39872 (90.0)'
$ws.Range('AF12').Value = '*maa://28932 (78.46), *maa://20106 (63.96), *maa://22769 (64.29)'
$ws.Range('D13').Value = 'maa://24999 (91.75), maa://36673 (92.65), maa://25001 (85.51)'
$ws.Range('L14').Value = 'maa://26245 (96.43), maa://21288 (96.27), maa://39841 (94.94), maa://36682 (97.37)'
$ws.Range('X14').Value = 'maa://37468 (90.0)'
$ws.Range('D15').Value = '*maa://22743 (76.96), maa://22734 (83.9), *maa://30808 (64.52), ***maa://36048 (28.57)'
$ws.Range('AF15').Value = 'maa://21364 (80.72), *maa://22766 (70.91), *maa://36666 (78.57)'
$ws.Range('D16').Value = 'maa://21441 (96.31), maa://36679 (93.02), maa://37650 (96.88)'
$ws.Range('AB16').Value = 'maa://26228 (95.4)'
$ws.Range('T17').Value = '**maa://42324 (46.15)'
$ws.Range('D18').Value = 'maa://24570 (97.07)'
$ws.Range('H18').Value = 'maa://24421 (90.3)'
$ws.Range('L18').Value = 'maa://22466 (88.73), *maa://22732 (50.6)'
$ws.Range('X18').Value = 'maa://21917 (97.8), maa://22741 (83.33)'
$ws.Range('AB19').Value = '*maa://30709 (62.96), *maa://36668 (55.84)'
$ws.Range('D20').Value = 'maa://21432 (89.86), maa://25198 (93.0), *maa://20795 (51.18), maa://36680 (96.55)'
$ws.Range('H20').Value = 'maa://22864 (89.12)'
$ws.Range('L20').Value = 'maa://41331 (84.21)'
$ws.Range('AB21').Value = '*maa://21443 (80.0), ***maa://23820 (29.82)'
$ws.Range('AF21').Value = 'maa://22524 (94.58), *maa://22432 (76.67)'
$ws.Range('L23').Value = 'maa://39756 (93.95), maa://39875 (93.65)'
$ws.Range('P23').Value = 'maa://30587 (91.94), *maa://29748 (75.59), ***maa://29785 (16.42), *maa://37566 (71.43)'
$ws.Range('D24').Value = '*maa://24368 (79.66)'
$ws.Range('X24').Value = 'maa://29988 (86.67), maa://23504 (93.23), **maa://22892 (39.58), *maa://25141 (77.6), maa://36663 (80.6), ***maa://22815 (23.08)'
$ws.Range('AF24').Value = 'maa://22523 (85.42), maa://36672 (80.77), maa://29910 (92.45), **maa://21440 (34.55)'
$ws.Range('T25').Value = 'maa://20109 (92.35), maa://22545 (100.0), maa://42915 (100.0)'
$ws.Range('X25').Value = '*maa://29890 (76.19)'
$ws.Range('AB25').Value = 'maa://31215 (85.57), *maa://24516 (79.78), maa://26001 (87.5)'
$ws.Range('D26').Value = 'maa://41802 (92.31)'
$ws.Range('AB26').Value = 'maa://42235 (93.15)'
$ws.Range('H27').Value = '**maa://21283 (48.0), maa://34494 (96.43), *maa://39601 (76.47), **maa://36665 (44.44)'
$ws.Range('T28').Value = 'maa://23263 (94.9), *maa://29765 (60.53)'
$ws.Range('X28').Value = 'maa://39929 (89.68), ***maa://39723 (14.29), maa://41749 (90.2)'
$ws.Range('AF28').Value = 'maa://36660 (92.88), *maa://36701 (64.29)'
$ws.Range('L29').Value = 'maa://28432 (92.81), *maa://28440 (76.6), maa://31400 (100.0), *maa://28650 (71.43)'
$ws.Range('AB30').Value = 'maa://42979 (96.23)'
$ws.Range('L31').Value = 'maa://35926 (93.94), maa://36258 (83.16), *maa://43904 (71.43)'
$ws.Range('T32').Value = 'maa://42859 (96.1), maa://41108 (87.76), maa://41238 (96.1)'
$ws.Range('T34').Value = 'maa://24526 (93.55)'
$ws.Range('L35').Value = 'maa://41296 (96.49)'
$ws.Range('AF38').Value = 'maa://36697 (86.19)'
$ws.Range('H39').Value = 'maa://25199 (84.82), maa://36670 (86.9), maa://30434 (89.06), ***maa://25036 (16.0), **maa://44165 (50.0)'
$ws.Range('P39').Value = 'maa://24709 (92.06)'
$ws.Range('P41').Value = '**maa://35616 (38.24), maa://43177 (84.62)'
$ws.Range('H46').Value = 'maa://35931 (92.47), maa://43901 (87.5)'
$ws.Range('H47').Value = 'maa://27410 (96.17), maa://29661 (97.83), maa://28038 (84.62)'
$ws.Range('H59').Value = 'maa://27746 (83.02), maa://31270 (94.78)'
$ws.Range('H60').Value = '*maa://40438 (58.14)'
$ws.Range('H62').Value = 'maa://42981 (96.15), maa://43903 (100.0)'
$ws.Range('H64').Value = 'maa://44405 (95.0)'
